# B1--and-B2-PowerPoint.pptx
#
# 1) Slide 5's table (the "Type of document" table) is switched from the
#    deck's custom table style to a built-in PowerPoint table style, via
#    the same Table.ApplyStyle(styleId) call that PowerPoint issues when a
#    user clicks a different swatch in the Table Design > Table Styles
#    gallery.
#
# 2) The presentation's applied Office Theme palette is changed back to
#    the default "Office" color scheme (it was the "Red Violet" variant of
#    the "Integral" theme). We drive this the supported way for this
#    object model: by writing each of the twelve theme colors in-order via
#    ThemeColorScheme.Colors(i).RGB (dk1, lt1, dk2, lt2, accent1-6, hlink,
#    folHlink) on the slide master's Theme -- the font scheme and format
#    scheme (fills/lines/effects) are identical between the two themes, so
#    nothing else needs to change there.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 -------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{37C7D9B7-4C43-4212-815C-74CA255E1328}")

# --- 2. Theme colors back to the default Office palette -----------------
function ConvertTo-OleRgb([string]$HexColor) {
    $r = [Convert]::ToInt32($HexColor.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($HexColor.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($HexColor.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ConvertTo-OleRgb $officeColors[$i - 1]
}
